# Update the DC-to-DC converter part numbers in the BOM, and related notes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Alimentation 5V part number change
$ws.Range("B10").Value = "PYBE20-Q24-S5-T"

# Row 11: Alimentation 24V part number change
$ws.Range("B11").Value = "PYBE20-Q24-S24-T"
$ws.Range("D11").Value = "Seulement pour le Sonar*2 "

# New footnote row for the 24V supply note
$ws.Range("A26").Value = "si oui l'alimentation peut etre retirer."

# Row 7: fix sonar model spelling (Imaginex 852 -> Imagenex 852)
$ws.Range("B7").Value = "Imagenex 852"

# New footnote row referencing the sonar power note
$ws.Range("A25").Value = "*2 a valider avec Imagenex si les nouveau sonar que vous avez peuvent fonctionner a 12v "

# Update the active selection to match the saved view state
$ws.Range("C5").Select()
